$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note cell below "Air_Location" header, added first so the shared-string
# table picks up "This HAS to be consistent" before "WaterT_C".
$ws.Range("I2").Value = "This HAS to be consistent"

# Header H1 was "Water_Temp_C" -> renamed to "WaterT_C"
$ws.Range("H1").Value = "WaterT_C"

# Column I widened to fit the new, longer text in I2 (was best-fit to
# "Air_Location"; now best-fit to "This HAS to be consistent").
$ws.Columns.Item(9).ColumnWidth = 20.75

# Selection moves to H9 to match the saved workbook state.
$ws.Range("H9").Select()
